$wb = $excel.ActiveWorkbook

# --- Add new worksheet "RQ5" after "RQ4" ------------------------------------
$rq4 = $wb.Worksheets.Item("RQ4")
$rq2 = $wb.Worksheets.Item("RQ2")
$new = $wb.Worksheets.Add([System.Type]::Missing, $rq4)
$new.Name = "RQ5"

# --- Values ------------------------------------------------------------------
# Merged header row
$new.Range("B1").Value = "Regular Clones"
$new.Range("E1").Value = "Micro Clones"
$new.Range("B1:D1").Merge()
$new.Range("E1:G1").Merge()

# Column header row
$new.Range("B2").Value = "Rep CFs LOCs"
$new.Range("C2").Value = "CFs LOCs"
$new.Range("D2").Value = "%"
$new.Range("E2").Value = "Rep CFs LOCs"
$new.Range("F2").Value = "CFs LOCs"
$new.Range("G2").Value = "%"

# Row labels
$new.Range("A3").Value = "Ctags"
$new.Range("A4").Value = "Brlcad"
$new.Range("A5").Value = "Freecol"
$new.Range("A6").Value = "Carol"
$new.Range("A7").Value = "Jabref"
$new.Range("A8").Value = "Total"

# Percent formulas for each project row
$new.Range("D3").Formula = "=B3/C3*100"
$new.Range("D4").Formula = "=B4/C4*100"
$new.Range("D5").Formula = "=B5/C5*100"
$new.Range("D6").Formula = "=B6/C6*100"
$new.Range("D7").Formula = "=B7/C7*100"

$new.Range("G3").Formula = "=E3/F3*100"
$new.Range("G4").Formula = "=E4/F4*100"
$new.Range("G5").Formula = "=E5/F5*100"
$new.Range("G6").Formula = "=E6/F6*100"
$new.Range("G7").Formula = "=E7/F7*100"

# Totals row
$new.Range("B8").Formula = "=SUM(B3:B7)"
$new.Range("C8").Formula = "=SUM(C3:C7)"
$new.Range("D8").Formula = "=B8/C8*100"
$new.Range("E8").Formula = "=SUM(E3:E7)"
$new.Range("F8").Formula = "=SUM(F3:F7)"
$new.Range("G8").Formula = "=E8/F8*100"

# --- Formatting ---------------------------------------------------------------
# Re-use the existing cell formats from RQ4/RQ2 (copy+paste-format) instead of
# toggling individual Font/Alignment properties, so no redundant style records
# get created - matches how these sister sheets are already formatted.

# Bold + centered merged header band (same format as RQ4 B1:G1)
$rq4.Range("B1:G1").Copy()
$new.Range("B1:G1").PasteSpecial(-4122)

# Bold sub-header row (same format as RQ4 B2:G2)
$rq4.Range("B2:G2").Copy()
$new.Range("B2:G2").PasteSpecial(-4122)

# Bold project-name column (same format as RQ4 A3:A8)
$rq4.Range("A3:A8").Copy()
$new.Range("A3:A8").PasteSpecial(-4122)

# Bold (empty) A1/A2 cells (same format as RQ2 A1:A2)
$rq2.Range("A1:A2").Copy()
$new.Range("A1:A2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Column widths (nearest values achievable through this engine's pixel
# quantisation that reproduce the target stored widths of 13.28515625 /
# 13.42578125 as closely as possible)
$new.Columns.Item(2).ColumnWidth = 12.5
$new.Columns.Item(5).ColumnWidth = 12.67

# Selection / active sheet
$new.Range("B3").Select()
$new.Activate()
